$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for the moving columns (Fecha, Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Origen, Precio $/Kg) per data row.
# The underlying weekly data set was re-sorted/shuffled across rows 2-36; the
# static descriptor columns (A,B,C,E,F,G,H,N,Q,R) are unchanged per row.
$rows = @(
    @{ Row=2; D=44847; I="Primera"; J=1110; K=1400; L=1500; M=1450; O="Provincia de Quillota"; P=1450 },
    @{ Row=3; D=44839; I="Primera"; J=3400; K=1400; L=1500; M=1447; O="Provincia de Quillota"; P=1447 },
    @{ Row=4; D=44846; I="Primera"; J=1000; K=1400; L=1450; M=1428; O="Provincia de Quillota"; P=1428 },
    @{ Row=5; D=44175; I="Primera"; J=1500; K=1300; L=1300; M=1300; O="Provincia de Quillota"; P=1300 },
    @{ Row=6; D=44175; I="Segunda"; J=1450; K=1000; L=1000; M=1000; O="Provincia de Quillota"; P=1000 },
    @{ Row=7; D=44838; I="Primera"; J=1480; K=1400; L=1500; M=1461; O="Provincia de Quillota"; P=1461 },
    @{ Row=8; D=44841; I="Primera"; J=1260; K=1400; L=1500; M=1454; O="Provincia de Quillota"; P=1454 },
    @{ Row=9; D=44162; I="Primera"; J=1200; K=1300; L=1300; M=1300; O="Provincia de Quillota"; P=1300 },
    @{ Row=10; D=44162; I="Segunda"; J=800; K=1000; L=1000; M=1000; O="Provincia de Quillota"; P=1000 },
    @{ Row=11; D=44837; I="Primera"; J=4400; K=1400; L=1500; M=1450; O="Provincia de Quillota"; P=1450 },
    @{ Row=12; D=44159; I="Primera"; J=1100; K=1300; L=1300; M=1300; O="Provincia de Quillota"; P=1300 },
    @{ Row=13; D=44159; I="Segunda"; J=800; K=1000; L=1000; M=1000; O="Provincia de Quillota"; P=1000 },
    @{ Row=14; D=44160; I="Primera"; J=750; K=1300; L=1300; M=1300; O="Provincia de Quillota"; P=1300 },
    @{ Row=15; D=44160; I="Segunda"; J=850; K=1000; L=1000; M=1000; O="Provincia de Quillota"; P=1000 },
    @{ Row=16; D=44165; I="Primera"; J=720; K=1200; L=1200; M=1200; O="Provincia de Quillota"; P=1200 },
    @{ Row=17; D=44165; I="Segunda"; J=750; K=1000; L=1000; M=1000; O="Provincia de Quillota"; P=1000 },
    @{ Row=18; D=44161; I="Primera"; J=1600; K=1300; L=1300; M=1300; O="Provincia de Quillota"; P=1300 },
    @{ Row=19; D=44161; I="Segunda"; J=1850; K=1000; L=1000; M=1000; O="Provincia de Quillota"; P=1000 },
    @{ Row=20; D=44176; I="Primera"; J=2500; K=1200; L=1300; M=1256; O="Provincia de Quillota"; P=1256 },
    @{ Row=21; D=44176; I="Segunda"; J=1500; K=1000; L=1000; M=1000; O="Provincia de Quillota"; P=1000 },
    @{ Row=22; D=44848; I="Primera"; J=1750; K=1400; L=1500; M=1449; O="Provincia de Quillota"; P=1449 },
    @{ Row=23; D=44172; I="Primera"; J=600; K=1300; L=1300; M=1300; O="Provincia de Quillota"; P=1300 },
    @{ Row=24; D=44172; I="Segunda"; J=550; K=1000; L=1000; M=1000; O="Provincia de Quillota"; P=1000 },
    @{ Row=25; D=44181; I="Primera"; J=1000; K=1300; L=1300; M=1300; O="Provincia de Quillota"; P=1300 },
    @{ Row=26; D=44181; I="Segunda"; J=900; K=900; L=900; M=900; O="Provincia de Quillota"; P=900 },
    @{ Row=27; D=44174; I="Primera"; J=2800; K=1200; L=1250; M=1221; O="Provincia de Quillota"; P=1221 },
    @{ Row=28; D=44174; I="Segunda"; J=1300; K=1000; L=1000; M=1000; O="Provincia de Quillota"; P=1000 },
    @{ Row=29; D=44179; I="Primera"; J=980; K=1200; L=1200; M=1200; O="Región Metropolitana"; P=1200 },
    @{ Row=30; D=44168; I="Primera"; J=1200; K=1300; L=1300; M=1300; O="Provincia de Quillota"; P=1300 },
    @{ Row=31; D=44168; I="Segunda"; J=850; K=1000; L=1000; M=1000; O="Provincia de Quillota"; P=1000 },
    @{ Row=32; D=44169; I="Primera"; J=950; K=1300; L=1300; M=1300; O="Provincia de Quillota"; P=1300 },
    @{ Row=33; D=44169; I="Segunda"; J=800; K=1000; L=1000; M=1000; O="Provincia de Quillota"; P=1000 },
    @{ Row=34; D=44167; I="Primera"; J=1430; K=1200; L=1300; M=1248; O="Provincia de Quillota"; P=1248 },
    @{ Row=35; D=44167; I="Segunda"; J=350; K=1000; L=1000; M=1000; O="Provincia de Quillota"; P=1000 },
    @{ Row=36; D=44845; I="Primera"; J=1000; K=1300; L=1500; M=1396; O="Provincia de Quillota"; P=1396 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.D    # D: Fecha
    $ws.Cells.Item($r.Row, 9).Value = $r.I    # I: Calidad
    $ws.Cells.Item($r.Row, 10).Value = $r.J   # J: Volumen
    $ws.Cells.Item($r.Row, 11).Value = $r.K   # K: Precio minimo
    $ws.Cells.Item($r.Row, 12).Value = $r.L   # L: Precio maximo
    $ws.Cells.Item($r.Row, 13).Value = $r.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r.Row, 15).Value = $r.O   # O: Origen
    $ws.Cells.Item($r.Row, 16).Value = $r.P   # P: Precio $/Kg
}

Write-Host "Applied weekly re-sort to rows 2-36"